# Fix ltdc and ram
# Rewrites the Typography table (sheet "Typography") and the Translation
# table (sheet "Translation") to match the updated font/typography setup
# and the refreshed set of translated texts.

$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------
# Typography sheet (Table7, header row 3, data starts row 4)
# Columns: B=Typography Name, C=Font, D=Size, E=Bpp, F=Fallback Character,
#          G=Wildcard Characters, H=Widget Wildcard Characters,
#          I=Wildcard Ranges, J=Ellipsis Character
# ---------------------------------------------------------------------
$typoRows = @(
    @("Default",        "verdana.ttf", 20,  4, "?", "", "", "",            ""),
    @("Large",           "verdana.ttf", 40,  4, "?", "", "", "",            ""),
    @("Small",           "verdana.ttf", 10,  4, "?", "", "", "",            ""),
    @("Typography_00",   "verdana.ttf", 20,  4, "?", "", "", "0-9,A-Z,a-z", ""),
    @("Typography_01",   "swisop3.ttf", 300, 4, "?", "", "", "0-9,A-Z,a-z", ""),
    @("Typography_02",   "swisop3.ttf", 40,  4, "?", "", "", "0-9,A-Z,a-z", ""),
    @("Typography_03",   "swisop3.ttf", 60,  4, "?", "", "", "0-9,A-Z,a-z", ""),
    @("Typography_04",   "swisop3.ttf", 50,  4, "?", "'' '", "", "0-9,A-Z,a-z", "")
)

$startRow = 4
for ($i = 0; $i -lt $typoRows.Count; $i++) {
    $r = $startRow + $i
    $row = $typoRows[$i]
    $wsTypo.Range("B$r").Value = $row[0]
    $wsTypo.Range("C$r").Value = $row[1]
    $wsTypo.Range("D$r").Value = $row[2]
    $wsTypo.Range("E$r").Value = $row[3]
    $wsTypo.Range("F$r").Value = $row[4]
    $wsTypo.Range("G$r").Value = $row[5]
    $wsTypo.Range("H$r").Value = $row[6]
    $wsTypo.Range("I$r").Value = $row[7]
    $wsTypo.Range("J$r").Value = $row[8]
}

# ---------------------------------------------------------------------
# Translation sheet, header row 3, data starts row 4
# Columns: B=TEXT ID, C=TYPOGRAPHY NAME, D=ALIGNMENT, E=DIRECTION, F=GB
# ---------------------------------------------------------------------
# Note: values that look like pure numbers (e.g. "99.9", "2") are prefixed
# with a leading apostrophe so Excel stores them as text (t="s"), matching
# the shared-string cells in the target workbook, instead of numeric cells.
$transRows = @(
    @("SingleUseId1",  "Typography_01", "Center", "LTR", "<value>"),
    @("SingleUseId2",  "Typography_03", "Center", "LTR", "<value>"),
    @("SingleUseId3",  "Typography_02", "Left",   "LTR", "BATT:  <value> V"),
    @("SingleUseId4",  "Typography_02", "Left",   "LTR", "TCS : <value>"),
    @("SingleUseId6",  "Typography_04", "Center", "LTR", "100C"),
    @("SingleUseId7",  "Typography_04", "Center", "LTR", "<value>"),
    @("ResourceId1",   "Default",       "Left",   "LTR", "TESTOWY"),
    @("SingleUseId8",  "Typography_02", "Left",   "LTR", "CLT:  <value>"),
    @("SingleUseId9",  "Typography_02", "Left",   "LTR", "'99.9"),
    @("SingleUseId11", "Typography_02", "Left",   "LTR", "'69.0"),
    @("SingleUseId12", "Typography_02", "Left",   "LTR", "'3"),
    @("SingleUseId13", "Typography_03", "Left",   "LTR", "'12345"),
    @("SingleUseId14", "Typography_01", "Left",   "LTR", "'2")
)

$startRow = 4
for ($i = 0; $i -lt $transRows.Count; $i++) {
    $r = $startRow + $i
    $row = $transRows[$i]
    $wsTrans.Range("B$r").Value = $row[0]
    $wsTrans.Range("C$r").Value = $row[1]
    $wsTrans.Range("D$r").Value = $row[2]
    $wsTrans.Range("E$r").Value = $row[3]
    $wsTrans.Range("F$r").Value = $row[4]
}
